$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D, E) contain numeric-looking text values
# (e.g. "241.98", "1.000", "  +0.16%  "). Force just the specific cells
# being updated to keep a Text number format so Excel does not silently
# convert them into real numbers/percentages (which would lose exact
# formatting, trailing zeros, or introduce floating point noise).
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "E29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.314.12'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.874.66'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '0.7107'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '241.98'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.07853'
$ws.Range("E8").Value = '  +1.98%  '
$ws.Range("D9").Value = '0.3126'
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").Value = '25.18'
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("D11").Value = '0.08390'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").Value = '1.865.41'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '5.242'
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = '0.7174'
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("D15").Value = '91.20'
$ws.Range("D16").Value = '6.209'
$ws.Range("E16").Value = '  +4.16%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.000008335'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '29.312.34'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '240.57'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '2.125.61'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '7.780'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '0.1594'
$ws.Range("E25").Value = '  -2.09%  '
$ws.Range("D26").Value = '9.052'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").Value = '162.54'
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").Value = '18.53'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").Value = '4.420'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '4.348'
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = '1.208'
$ws.Range("E32").Value = '  -5.87%  '
$ws.Range("D33").Value = '0.05358'
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("D34").Value = '1.946'
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").Value = '1.177'
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").Value = '0.7481'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").Value = '1.291.37'
$ws.Range("E38").Value = '  +11.86%  '
$ws.Range("D39").Value = '0.01885'
$ws.Range("E39").Value = '  +1.60%  '
$ws.Range("D40").Value = '2.740'
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D41").Value = '6.569'
$ws.Range("E41").Value = '  +3.21%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '111.00'
$ws.Range("E42").Value = '  +5.57%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8945'
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").Value = '73.06'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '0.00000000131'
$ws.Range("E45").Value = '  +9.52%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = '2.023.53'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '1.800'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = '0.5196'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '9.454'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").Value = '0.4356'
$ws.Range("E51").Value = '  +1.38%  '
